$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17 currently is an empty row (matching the formatting of row 16/others).
# Copy the number/date formats from row 15 (A:G), which has the same layout
# (No / Date / Application / Task / % / Status / Comments), onto row 17 so the
# new entries pick up the correct date and percentage formatting.
$ws.Range("A15:G15").Copy()
$ws.Range("A17:G17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new daily-track entry for row 17.
$ws.Cells.Item(17, 1).Value = 9
$ws.Cells.Item(17, 2).Value = 44574
$ws.Cells.Item(17, 3).Value = "GSS"
$ws.Cells.Item(17, 4).Value = "1. Uploading PDF files task has been completed, tested and it is uploading smoothly"
$ws.Cells.Item(17, 5).Value = 1
$ws.Cells.Item(17, 6).Value = "Completed"
$ws.Cells.Item(17, 7).Value = $null

# Move the active selection to F17, matching the saved view state.
$null = $ws.Range("F17").Select()
